# Contamination / metals results workbook rework:
#   - The original Sheet1 (ESTUARY-level Cadmium/Cuivre/Mercure/Plomb medians,
#     with "median_ng_gww"/"median_1"/"median_2" columns) is duplicated into a
#     new "Sheet2" placed right after Sheet1. On that copy, the C/D header
#     labels are renamed to "First 5 years" / "Last 5 years", the old E
#     ("median_2") column is dropped, the old C ("median_ng_gww") values slide
#     into D ("Last 5 years") unchanged, and a new C ("First 5 years") value is
#     filled in per row.
#   - Sheet1 itself is cleared and replaced by a new, condensed per-parameter
#     summary table (PARAMETRE_LIBELLE / First 5 years / Last 5 years), one
#     row per metal, aggregated across estuaries.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Duplicate Sheet1 -> Sheet2 (keeps original data + header styling) ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# --- Sheet2: rename the C/D headers, drop column E ---
$ws2.Range("C1").Value = "First 5 years"
$ws2.Range("D1").Value = "Last 5 years"
$ws2.Columns.Item(5).Delete()

# --- Sheet2: new "First 5 years" (C) values; "Last 5 years" (D) already
#     holds the original median_ng_gww figures and needs no change ---
$sheet2FirstValues = @{
    2  = 7.7195
    3  = 0.2785
    4  = 0.6785
    5  = 4.491
    6  = 1.517
    7  = 1.774
    8  = 0.016
    9  = 0.019
    10 = 0.0495
    11 = 0.323
    12 = 0.866
    13 = 0.6335
}
foreach ($r in $sheet2FirstValues.Keys) {
    $ws2.Cells.Item($r, 3).Value = $sheet2FirstValues[$r]
}

# --- Sheet1: cleared and rebuilt as the condensed per-metal summary ---
$ws1.Cells.Clear()

$ws1.Range("A1").Value = "PARAMETRE_LIBELLE"
$ws1.Range("B1").Value = "First 5 years"
$ws1.Range("C1").Value = "Last 5 years"
$ws1.Range("A1:C1").Font.Bold = $true
$ws1.Range("A1:C1").HorizontalAlignment = -4108

$sheet1Data = @(
    @("Cadmium", 0.911,  0.23),
    @("Cuivre",  2.079,  2.069),
    @("Mercure", 0.019,  0.025),
    @("Plomb",   0.5875, 0.249)
)

$r = 2
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws1.Activate()
